$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStrings = @(
    'Tử Vi tọa thủ tại cung đối Quan Lộc',
    'Thiên Cơ tọa thủ tại cung đối Quan Lộc',
    'Thái Dương tọa thủ tại cung đối Quan Lộc',
    'Vũ Khúc tọa thủ tại cung đối Quan Lộc',
    'Thiên Đồng tọa thủ tại cung đối Quan Lộc',
    'Liêm Trinh tọa thủ tại cung đối Quan Lộc',
    'Thiên Phủ tọa thủ tại cung đối Quan Lộc',
    'Thái Âm tọa thủ tại cung đối Quan Lộc',
    'Tham Lang tọa thủ tại cung đối Quan Lộc',
    'Cự Môn tọa thủ tại cung đối Quan Lộc',
    'Thiên Tướng tọa thủ tại cung đối Quan Lộc',
    'Thiên Lương tọa thủ tại cung đối Quan Lộc',
    'Thất Sát tọa thủ tại cung đối Quan Lộc',
    'Phá Quân tọa thủ tại cung đối Quan Lộc',
    'Tử Vi đồng cung Thiên Cơ tại cung đối Quan Lộc',
    'Tử Vi đồng cung Thái Dương tại cung đối Quan Lộc',
    'Tử Vi đồng cung Vũ Khúc tại cung đối Quan Lộc',
    'Tử Vi đồng cung Thiên Đồng tại cung đối Quan Lộc',
    'Tử Vi đồng cung Liêm Trinh tại cung đối Quan Lộc',
    'Tử Vi đồng cung Thiên Phủ tại cung đối Quan Lộc',
    'Tử Vi đồng cung Thái Âm tại cung đối Quan Lộc',
    'Tử Vi đồng cung Tham Lang tại cung đối Quan Lộc',
    'Tử Vi đồng cung Cự Môn tại cung đối Quan Lộc',
    'Tử Vi đồng cung Thiên Tướng tại cung đối Quan Lộc',
    'Tử Vi đồng cung Thiên Lương tại cung đối Quan Lộc',
    'Tử Vi đồng cung Thất Sát tại cung đối Quan Lộc',
    'Tử Vi đồng cung Phá Quân tại cung đối Quan Lộc',
    'Thiên Cơ đồng cung Thái Dương tại cung đối Quan Lộc',
    'Thiên Cơ đồng cung Vũ Khúc tại cung đối Quan Lộc',
    'Thiên Cơ đồng cung Thiên Đồng tại cung đối Quan Lộc',
    'Thiên Cơ đồng cung Liêm Trinh tại cung đối Quan Lộc',
    'Thiên Cơ đồng cung Thiên Phủ tại cung đối Quan Lộc',
    'Thiên Cơ đồng cung Thái Âm tại cung đối Quan Lộc',
    'Thiên Cơ đồng cung Tham Lang tại cung đối Quan Lộc',
    'Thiên Cơ đồng cung Cự Môn tại cung đối Quan Lộc',
    'Thiên Cơ đồng cung Thiên Tướng tại cung đối Quan Lộc',
    'Thiên Cơ đồng cung Thiên Lương tại cung đối Quan Lộc',
    'Thiên Cơ đồng cung Thất Sát tại cung đối Quan Lộc',
    'Thiên Cơ đồng cung Phá Quân tại cung đối Quan Lộc',
    'Thái Dương đồng cung Vũ Khúc tại cung đối Quan Lộc',
    'Thái Dương đồng cung Thiên Đồng tại cung đối Quan Lộc',
    'Thái Dương đồng cung Liêm Trinh tại cung đối Quan Lộc',
    'Thái Dương đồng cung Thiên Phủ tại cung đối Quan Lộc',
    'Thái Dương đồng cung Thái Âm tại cung đối Quan Lộc',
    'Thái Dương đồng cung Tham Lang tại cung đối Quan Lộc',
    'Thái Dương đồng cung Cự Môn tại cung đối Quan Lộc',
    'Thái Dương đồng cung Thiên Tướng tại cung đối Quan Lộc',
    'Thái Dương đồng cung Thiên Lương tại cung đối Quan Lộc',
    'Thái Dương đồng cung Thất Sát tại cung đối Quan Lộc',
    'Thái Dương đồng cung Phá Quân tại cung đối Quan Lộc',
    'Vũ Khúc đồng cung Thiên Đồng tại cung đối Quan Lộc',
    'Vũ Khúc đồng cung Liêm Trinh tại cung đối Quan Lộc',
    'Vũ Khúc đồng cung Thiên Phủ tại cung đối Quan Lộc',
    'Vũ Khúc đồng cung Thái Âm tại cung đối Quan Lộc',
    'Vũ Khúc đồng cung Tham Lang tại cung đối Quan Lộc',
    'Vũ Khúc đồng cung Cự Môn tại cung đối Quan Lộc',
    'Vũ Khúc đồng cung Thiên Tướng tại cung đối Quan Lộc',
    'Vũ Khúc đồng cung Thiên Lương tại cung đối Quan Lộc',
    'Vũ Khúc đồng cung Thất Sát tại cung đối Quan Lộc',
    'Vũ Khúc đồng cung Phá Quân tại cung đối Quan Lộc',
    'Thiên Đồng đồng cung Liêm Trinh tại cung đối Quan Lộc',
    'Thiên Đồng đồng cung Thiên Phủ tại cung đối Quan Lộc',
    'Thiên Đồng đồng cung Thái Âm tại cung đối Quan Lộc',
    'Thiên Đồng đồng cung Tham Lang tại cung đối Quan Lộc',
    'Thiên Đồng đồng cung Cự Môn tại cung đối Quan Lộc',
    'Thiên Đồng đồng cung Thiên Tướng tại cung đối Quan Lộc',
    'Thiên Đồng đồng cung Thiên Lương tại cung đối Quan Lộc',
    'Thiên Đồng đồng cung Thất Sát tại cung đối Quan Lộc',
    'Thiên Đồng đồng cung Phá Quân tại cung đối Quan Lộc',
    'Liêm Trinh đồng cung Thiên Phủ tại cung đối Quan Lộc',
    'Liêm Trinh đồng cung Thái Âm tại cung đối Quan Lộc',
    'Liêm Trinh đồng cung Tham Lang tại cung đối Quan Lộc',
    'Liêm Trinh đồng cung Cự Môn tại cung đối Quan Lộc',
    'Liêm Trinh đồng cung Thiên Tướng tại cung đối Quan Lộc',
    'Liêm Trinh đồng cung Thiên Lương tại cung đối Quan Lộc',
    'Liêm Trinh đồng cung Thất Sát tại cung đối Quan Lộc',
    'Liêm Trinh đồng cung Phá Quân tại cung đối Quan Lộc',
    'Thiên Phủ đồng cung Thái Âm tại cung đối Quan Lộc',
    'Thiên Phủ đồng cung Tham Lang tại cung đối Quan Lộc',
    'Thiên Phủ đồng cung Cự Môn tại cung đối Quan Lộc',
    'Thiên Phủ đồng cung Thiên Tướng tại cung đối Quan Lộc',
    'Thiên Phủ đồng cung Thiên Lương tại cung đối Quan Lộc',
    'Thiên Phủ đồng cung Thất Sát tại cung đối Quan Lộc',
    'Thiên Phủ đồng cung Phá Quân tại cung đối Quan Lộc',
    'Thái Âm đồng cung Tham Lang tại cung đối Quan Lộc',
    'Thái Âm đồng cung Cự Môn tại cung đối Quan Lộc',
    'Thái Âm đồng cung Thiên Tướng tại cung đối Quan Lộc',
    'Thái Âm đồng cung Thiên Lương tại cung đối Quan Lộc',
    'Thái Âm đồng cung Thất Sát tại cung đối Quan Lộc',
    'Thái Âm đồng cung Phá Quân tại cung đối Quan Lộc',
    'Tham Lang đồng cung Cự Môn tại cung đối Quan Lộc',
    'Tham Lang đồng cung Thiên Tướng tại cung đối Quan Lộc',
    'Tham Lang đồng cung Thiên Lương tại cung đối Quan Lộc',
    'Tham Lang đồng cung Thất Sát tại cung đối Quan Lộc',
    'Tham Lang đồng cung Phá Quân tại cung đối Quan Lộc',
    'Cự Môn đồng cung Thiên Tướng tại cung đối Quan Lộc',
    'Cự Môn đồng cung Thiên Lương tại cung đối Quan Lộc',
    'Cự Môn đồng cung Thất Sát tại cung đối Quan Lộc',
    'Cự Môn đồng cung Phá Quân tại cung đối Quan Lộc',
    'Thiên Tướng đồng cung Thiên Lương tại cung đối Quan Lộc',
    'Thiên Tướng đồng cung Thất Sát tại cung đối Quan Lộc',
    'Thiên Tướng đồng cung Phá Quân tại cung đối Quan Lộc',
    'Thiên Lương đồng cung Thất Sát tại cung đối Quan Lộc',
    'Thiên Lương đồng cung Phá Quân tại cung đối Quan Lộc',
    'Thất Sát đồng cung Phá Quân tại cung đối Quan Lộc'
)

$n = $newStrings.Length
$arr = New-Object 'object[,]' $n,2
for ($i = 0; $i -lt $n; $i++) {
    $arr[$i,0] = $newStrings[$i]
    $arr[$i,1] = $newStrings[$i]
}

# First block: rows 4304-4317 (14 single-star entries), row 4303 left blank as a separator
$singleCount = 14
$singleArr = New-Object 'object[,]' $singleCount,2
for ($i = 0; $i -lt $singleCount; $i++) {
    $singleArr[$i,0] = $arr[$i,0]
    $singleArr[$i,1] = $arr[$i,1]
}
$ws.Range("A4304:B4317").Value = $singleArr

# Second block: rows 4319-4409 (91 pair entries), row 4318 left blank as a separator
$pairCount = 91
$pairArr = New-Object 'object[,]' $pairCount,2
for ($i = 0; $i -lt $pairCount; $i++) {
    $pairArr[$i,0] = $arr[$singleCount + $i,0]
    $pairArr[$i,1] = $arr[$singleCount + $i,1]
}
$ws.Range("A4319:B4409").Value = $pairArr

# Match the author's final selection/view state captured in the sheet
$win = $excel.ActiveWindow
$win.ScrollRow = 4300
$win.ScrollColumn = 1
$ws.Range("B4319:B4409").Select() | Out-Null

